# Publication details update: JUL update, lab foundation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Total papers" count: 21 -> 22
$ws.Range("B2").Value = 22

# Update "Senior, Co-senior or primary supervisor" count: 2 -> 4
$ws.Range("B4").Value = 4

# Rename "As corresponding author" -> "As (co-)corresponding author"
$ws.Range("A7").Value = "As (co-)corresponding author"

# Remove the now-obsolete "Papers in review as Co-First/Final author" row (row 8)
$ws.Rows.Item(8).Delete()

# Fix up the active selection to match the saved state
$ws.Range("B8").Select()
